# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7625
$ws1.Range("F4").Value = 219
$ws1.Range("F5").Value = 30
$ws1.Range("F6").Value = 273
$ws1.Range("F7").Value = 1145
$ws1.Range("F8").Value = 201
$ws1.Range("F10").Value = 150
$ws1.Range("F11").Value = 38

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7625
$ws4.Range("F4").Value = 219
$ws4.Range("F5").Value = 30
$ws4.Range("F6").Value = 273
$ws4.Range("F7").Value = 1145
$ws4.Range("F8").Value = 201
$ws4.Range("F11").Value = 150
$ws4.Range("F12").Value = 38
